$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "NIG(0.9004031587145218, 0.6405938341949465, 1.3646639060473915, 2.9451809548732046)"
$ws.Range("C2").Value = "JSU(-0.8869957803715276, 0.9937017812981697, 4.45740534070488, 4.347987908665581)"
$ws.Range("D2").Value = "JSU(-0.6753344863133022, 0.9784696193265736, 1.0577739163677964, 2.2720353053163134)"
$ws.Range("E2").Value = "NIG(0.9418662832561067, 0.7438333367768206, 4.36951242257131, 5.346253263229717)"
